# Generate Report for Handback
#
# Marks the a.md / b.md rows as handed back (in sync with en-US) across the
# Overview sheet and the per-locale (zh-cn / de-de) sheets, and records the
# handback target file / handback file / handback datetime for the "a.md"
# source on each locale sheet (b.md was a duplicate of a.md's content, so it
# shares the same handback target).

$wb = $excel.ActiveWorkbook

$statusHandedBack = "Handed back: in sync with en-US"
$repoBlobBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/81989117afd2eca6ccaa77c0ebac7f7f34eef237/e2e"

# ---------------------------------------------------------------------
# Overview sheet: zh-cn / de-de status columns for both rows
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $statusHandedBack
$wsOverview.Range("F2").Value = $statusHandedBack
$wsOverview.Range("E3").Value = $statusHandedBack
$wsOverview.Range("F3").Value = $statusHandedBack
$wsOverview.Columns.Item(5).ColumnWidth = 29.9777047293527 - 0.8333333333333333
$wsOverview.Columns.Item(6).ColumnWidth = 29.9777047293527 - 0.8333333333333333

function Update-LocaleSheet {
    param(
        [string]$SheetName,
        [string]$XliffFileName,
        [string]$HandbackDateTime
    )

    $ws = $wb.Worksheets.Item($SheetName)

    # Status column for both data rows
    $ws.Range("C2").Value = $statusHandedBack
    $ws.Range("C3").Value = $statusHandedBack

    # Latest Handback File / Latest Handback DateTime
    $ws.Range("J2").Value = $XliffFileName
    $ws.Range("J3").Value = $XliffFileName
    $ws.Range("K2").Value = $HandbackDateTime
    $ws.Range("K3").Value = $HandbackDateTime

    # Widen the Status and Latest Handback File columns to fit the new text
    $ws.Columns.Item(3).ColumnWidth = 29.9777047293527 - 0.8333333333333333
    $ws.Columns.Item(10).ColumnWidth = 40 - 0.8333333333333333

    # Re-create every hyperlink on the sheet, in row/column order, so the
    # new "Latest Target File" links (I2/I3, pointing at a.md) line up
    # alongside the existing Source File Name links (A2/A3).
    $ws.Hyperlinks.Delete()
    $ws.Hyperlinks.Add($ws.Range("A2"), "$repoBlobBase/a.md", "", "", "a.md")
    $ws.Hyperlinks.Add($ws.Range("I2"), "$repoBlobBase/a.md", "", "", "a.md")
    $ws.Hyperlinks.Add($ws.Range("A3"), "$repoBlobBase/b.md", "", "", "b.md")
    $ws.Hyperlinks.Add($ws.Range("I3"), "$repoBlobBase/a.md", "", "", "a.md")

    # I2/I3 now hold the "Latest Target File" name (a.md)
    $ws.Range("I2").Value = "a.md"
    $ws.Range("I3").Value = "a.md"
}

Update-LocaleSheet "zh-cn" "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf" "2016-08-18 22:36:40"
Update-LocaleSheet "de-de" "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf" "2016-08-18 22:36:47"
